$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "65.883.28"
$ws.Range("E2").Value = "  -2.34%  "

# Row 3
$ws.Range("D3").Value = "3.476.87"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "582.08"
$ws.Range("E5").Value = "  -1.53%  "

# Row 6
Set-TextValue $ws.Range("D6") "171.75"
$ws.Range("E6").Value = "  -3.85%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.597"
$ws.Range("E8").Value = "  -1.24%  "

# Row 9
$ws.Range("D9").Value = "3.472.32"
$ws.Range("E9").Value = "  +0.91%  "

# Row 10
$ws.Range("E10").Value = "  -5.78%  "

# Row 11
$ws.Range("E11").Value = "  -1.65%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.410"
$ws.Range("E12").Value = "  -3.70%  "

# Row 13
$ws.Range("D13").Value = "4.076.17"
$ws.Range("E13").Value = "  +0.77%  "

# Row 14
$ws.Range("E14").Value = "  +1.23%  "

# Row 15
Set-TextValue $ws.Range("D15") "29.78"
$ws.Range("E15").Value = "  -6.67%  "

# Row 16
$ws.Range("D16").Value = "65.920.58"
$ws.Range("E16").Value = "  -2.28%  "

# Row 17
$ws.Range("E17").Value = "  -3.28%  "

# Row 18
$ws.Range("D18").Value = "3.474.36"
$ws.Range("E18").Value = "  +0.93%  "

# Row 19
$ws.Range("E19").Value = "  -3.30%  "

# Row 20
Set-TextValue $ws.Range("D20") "13.86"
$ws.Range("E20").Value = "  -0.80%  "

# Row 21
Set-TextValue $ws.Range("D21") "366.00"
$ws.Range("E21").Value = "  -5.01%  "

# Row 22
Set-TextValue $ws.Range("D22") "7.72"
$ws.Range("E22").Value = "  -1.24%  "

# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D24") "72.50"
$ws.Range("E24").Value = "  +1.99%  "

# Row 25
$ws.Range("E25").Value = "  +5.75%  "

# Row 26
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.55"
$ws.Range("E27").Value = "  -6.27%  "

# Row 28
$ws.Range("E28").Value = "  +2.53%  "

# Row 29
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
Set-TextValue $ws.Range("D30") "23.99"
$ws.Range("E30").Value = "  +2.41%  "

# Row 31
$ws.Range("E31").Value = "  -4.86%  "

# Row 32
$ws.Range("E32").Value = "  -2.97%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("E34").Value = "  -1.09%  "

# Row 35
$ws.Range("E35").Value = "  -6.03%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.52"
$ws.Range("E36").Value = "  -1.53%  "

# Row 37
Set-TextValue $ws.Range("D37") "160.48"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38
Set-TextValue $ws.Range("D38") "29.30"
$ws.Range("E38").Value = "  +13.84%  "

# Row 39
$ws.Range("E39").Value = "  +1.02%  "

# Row 40
$ws.Range("D40").Value = "2.821.27"
$ws.Range("E40").Value = "  +4.62%  "

# Row 41
$ws.Range("E41").Value = "  -5.27%  "

# Row 42
Set-TextValue $ws.Range("D42") "6.47"
$ws.Range("E42").Value = "  -2.03%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.55"
$ws.Range("E43").Value = "  -6.69%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "4.44"
$ws.Range("E44").Value = "  -1.77%  "

# Row 45
$ws.Range("E45").Value = "  -4.37%  "

# Row 46
Set-TextValue $ws.Range("D46") "40.07"
$ws.Range("E46").Value = "  -2.60%  "

# Row 47
Set-TextValue $ws.Range("D47") "24.03"
$ws.Range("E47").Value = "  -7.07%  "

# Row 48
$ws.Range("E48").Value = "  -2.85%  "

# Row 49
Set-TextValue $ws.Range("D49") "324.57"
$ws.Range("E49").Value = "  -0.31%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.813"
$ws.Range("E50").Value = "  -2.53%  "

# Row 51
$ws.Range("E51").Value = "  -2.69%  "
